# Rename the two "dFBA net ..." sheets to "dFBA objective ..." per the
# wc_lang DfbaNetReaction -> DfbaObjReaction / DfbaNetSpecies -> DfbaObjSpecies
# rename described in the commit message.
$wb = $excel.ActiveWorkbook

$reactionsSheet = $wb.Worksheets.Item("dFBA net reactions")
$reactionsSheet.Name = "dFBA objective reactions"

$speciesSheet = $wb.Worksheets.Item("dFBA net species")
$speciesSheet.Name = "dFBA objective species"

# The column header on the "species" sheet still reads the old label -
# update the cell text to match the new terminology.
$speciesSheet.Range("C1").Value = "dFBA objective reaction"

# Make the renamed species sheet the active tab/selection (this was the
# last sheet the author was looking at when the workbook was saved).
$speciesSheet.Select()
$speciesSheet.Range("F8").Select()

# Resize the saved window to match the author's larger Excel window.
$win = $excel.ActiveWindow
$win.Width = 28200
$win.Height = 14070
